# Add f6/f7 pass & p-value results for row 13 (rhyme (hashtag separated) row)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# f6 pass / f6 p value
$ws.Range("L13").Value = 1
$ws.Range("L13").Style = "Good"

$ws.Range("M13").Value = 0.00010011
$ws.Range("M13").NumberFormat = "0.00E+00"

# f7 pass / f7 p value
$ws.Range("N13").Value = 1
$ws.Range("N13").Style = "Good"
$ws.Range("N13").NumberFormat = "0.00E+00"

$ws.Range("O13").Value = 0.0035

# update the selected cell shown in the sheet view
$ws.Range("N10").Select()
